$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (Journal_Entry grouping id) for rows 11-101
$ws.Range("B11:B18").Value = 1510
$ws.Range("B19:B29").Value = 2010
$ws.Range("B30:B35").Value = 2200
$ws.Range("B36:B45").Value = 3010
$ws.Range("B46:B55").Value = 3110
$ws.Range("B56:B74").Value = 4010
$ws.Range("B75:B84").Value = 5010
$ws.Range("B85:B101").Value = 6010

# Update the view: top-left cell and selected cell/range
$excel.ActiveWindow.ScrollRow = 68
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F89").Select()
